$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCall")

# --- Capture state that will shift when we insert the two new columns ---
# The "custom fields" comment currently lives on L1; once we insert two
# columns in front of it, it will be pushed out to N1 (but the comment's
# anchor itself does not move automatically, so we recreate it later).
$existingComment = $ws.Range("L1").Comment
$commentText = $existingComment.Text()

# Width used by column K - the two new columns should match it.
$callBasisColumnWidth = $ws.Range("K1").EntireColumn.ColumnWidth

# --- Insert two new columns (for "Send Payment Notification" / "Send Call Notice") ---
$ws.Range("L1:M1").EntireColumn.Insert()

# Give the two new columns the same width as the neighbouring "Call Basis" column.
$ws.Range("L1").EntireColumn.ColumnWidth = $callBasisColumnWidth
$ws.Range("M1").EntireColumn.ColumnWidth = $callBasisColumnWidth

# Move the "custom fields" comment from its old spot (now shifted to N1) onto the
# new, empty N1 cell - recreate it there with the same text.
$ws.Range("L1").Comment.Delete()
$ws.Range("N1").AddComment($commentText)

# --- New header row values ---
$ws.Range("L1").Value = "Send Payment Notification"
$ws.Range("M1").Value = "Send Call Notice"

# --- New column values for the existing data rows ---
$ws.Range("L2").Value = "Yes"
$ws.Range("M2").Value = "Yes"
$ws.Range("L3").Value = "Yes"
$ws.Range("M3").Value = "Yes"
$ws.Range("L4").Value = "Yes"
$ws.Range("M4").Value = "Yes"

# --- Extend the "Percentage of Commitment,Upload" list validation to cover the new columns ---
$ws.Range("J2:J1048576").Validation.Delete()
$ws.Range("K6:K1048576").Validation.Delete()
$ws.Range("J2:J1048576").Validation.Add(3, 1, 1, """Percentage of Commitment,Upload""")
$ws.Range("K6:M1048576").Validation.Add(3, 1, 1, """Percentage of Commitment,Upload""")

# --- Match the saved selection state ---
$ws.Range("A4").Select()
